$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the "Serie" date column) and append
# the new quarterly record right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$cellA = $ws.Cells.Item($newRow, 1)

# A plain Value assignment of "01-07-2021" would be auto-recognized by
# Excel as a date and stored as a serial number (plus a new number-format
# style), but the source data keeps these quarter labels as plain text in
# the shared-string table. Entering it as a formula that evaluates to the
# literal text sidesteps the date auto-detection, then Copy/PasteSpecial
# values collapses the formula back down to a plain text constant, so the
# cell ends up as an ordinary (unstyled) shared string, exactly like the
# rows above it.
$cellA.Formula = "=""01-07-2021"""
$cellA.Copy()
$cellA.PasteSpecial(-4163)  # xlPasteValues
$ws.Application.CutCopyMode = $false

$ws.Cells.Item($newRow, 2).Value = 689
